# Generate Report for Handback
# The 827e4b46-... file has now been handed back (in sync with en-US)
# for both the zh-cn and de-de locales. Update the per-locale detail
# sheets' Status + Latest Handback DateTime, and roll the change up
# into the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 is the 827e4b46-...md entry ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("G2").Value = "2016-03-09 13:28:01"
$zh.Range("G3").Value = "2016-03-09 13:28:01"

# --- de-de sheet: row 3 is the 827e4b46-...md entry ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("G2").Value = "2016-03-09 13:28:18"
$de.Range("G3").Value = "2016-03-09 13:28:18"

# --- Overview sheet: roll the per-locale status up ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"
